# "final re-run including fig: unmet social support"
#
# The chart on slide 1 (results/figures/pptx/impact_personal_relations.pptx)
# is a bar chart built from individual drawn shapes grouped together.
# The underlying counts shifted slightly on re-run (N=146 -> N=143, and a
# few per-bar counts/percentages changed by 1), which reflows the gridlines,
# bars, data-label positions/sizes and axis tick labels. This script pokes
# each affected shape (accessed via the single top-level group's
# GroupItems collection) to match the new layout/values.
#
# Note: shape Top/Height are expressed in points through the COM object
# model (1 pt = 12700 EMU) and are rounded through a 32-bit float
# internally, so the literal point values below were chosen (via a small
# offline search) so that they land exactly on the target EMU value after
# that rounding, instead of using a naive value/12700 division.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$grp = $s.Shapes.Item(1)
$grp.GroupItems.Item(3).Top = 394.855224609375  # pl5 y 5036889 -> 5014661
$grp.GroupItems.Item(4).Top = 310.27520751953125  # pl6 y 3984952 -> 3940495
$grp.GroupItems.Item(5).Top = 225.69520568847656  # pl7 y 2933014 -> 2866329
$grp.GroupItems.Item(6).Top = 141.1152801513672  # pl8 y 1881077 -> 1792164
$grp.GroupItems.Item(7).Top = 461.6911315917969  # rc9 y 5944726 -> 5863477
$grp.GroupItems.Item(7).Height = 17.744016647338867  # rc9 cy 144101 -> 225349
$grp.GroupItems.Item(8).Top = 384.8002624511719  # rc10 y 4863968 -> 4886963
$grp.GroupItems.Item(8).Height = 94.6348876953125  # rc10 cy 1224858 -> 1201863
$grp.GroupItems.Item(9).Top = 367.05615234375  # rc11 y 4647816 -> 4661613
$grp.GroupItems.Item(9).Height = 112.37898254394531  # rc11 cy 1441010 -> 1427213
$grp.GroupItems.Item(10).Top = 278.3359069824219  # rc12 y 3567059 -> 3534866
$grp.GroupItems.Item(10).Height = 201.09921264648438  # rc12 cy 2521767 -> 2553960
$grp.GroupItems.Item(12).Top = 337.4827575683594  # rc14 y 4431665 -> 4286031
$grp.GroupItems.Item(12).Height = 141.95237731933594  # rc14 cy 1657161 -> 1802795
$grp.GroupItems.Item(13).Top = 415.4645690917969  # tx15 y 5359942 -> 5276400
$grp.GroupItems.Item(13).Height = 10.40574836730957  # tx15 cy 129860 -> 132153
$grp.GroupItems.Item(13).TextFrame.TextRange.Text = "3"  # tx15 text '2' -> '3'
$grp.GroupItems.Item(14).Top = 433.0052185058594  # tx16 y 5580415 -> 5499166
$grp.GroupItems.Item(14).TextFrame.TextRange.Text = "(2%)"  # tx16 text '(1%)' -> '(2%)'
$grp.GroupItems.Item(15).Top = 338.5806579589844  # tx17 y 4279184 -> 4299974
$grp.GroupItems.Item(15).Height = 10.398818969726562  # tx17 cy 129860 -> 132065
$grp.GroupItems.Item(15).TextFrame.TextRange.Text = "16"  # tx17 text '17' -> '16'
$grp.GroupItems.Item(16).Top = 356.1143493652344  # tx18 y 4499657 -> 4522652
$grp.GroupItems.Item(16).TextFrame.TextRange.Text = "(11%)"  # tx18 text '(12%)' -> '(11%)'
$grp.GroupItems.Item(17).Top = 320.8365478515625  # tx19 y 4060827 -> 4074624
$grp.GroupItems.Item(17).TextFrame.TextRange.Text = "19"  # tx19 text '20' -> '19'
$grp.GroupItems.Item(18).Top = 338.3703308105469  # tx20 y 4283506 -> 4297303
$grp.GroupItems.Item(18).TextFrame.TextRange.Text = "(13%)"  # tx20 text '(14%)' -> '(13%)'
$grp.GroupItems.Item(19).Top = 232.109375  # tx21 y 2979981 -> 2947789
$grp.GroupItems.Item(19).TextFrame.TextRange.Text = "34"  # tx21 text '35' -> '34'
$grp.GroupItems.Item(20).Top = 249.65000915527344  # tx22 y 3202748 -> 3170555
$grp.GroupItems.Item(21).Top = 155.4407196044922  # tx23 y 1971363 -> 1974097
$grp.GroupItems.Item(21).Height = 10.183465003967285  # tx23 cy 132065 -> 129330
$grp.GroupItems.Item(21).TextFrame.TextRange.Text = "47"  # tx23 text '49' -> '47'
$grp.GroupItems.Item(22).TextFrame.TextRange.Text = "(33%)"  # tx24 text '(34%)' -> '(33%)'
$grp.GroupItems.Item(23).Top = 291.4367980957031  # tx25 y 3844588 -> 3701247
$grp.GroupItems.Item(23).Height = 10.225197792053223  # tx25 cy 132153 -> 129860
$grp.GroupItems.Item(23).TextFrame.TextRange.Text = "24"  # tx25 text '23' -> '24'
$grp.GroupItems.Item(24).Top = 308.796875  # tx26 y 4067354 -> 3921720
$grp.GroupItems.Item(24).TextFrame.TextRange.Text = "(17%)"  # tx26 text '(16%)' -> '(17%)'
$grp.GroupItems.Item(27).Top = 390.37860107421875  # tx29 y 4980037 -> 4957808
$grp.GroupItems.Item(28).Top = 305.7986755371094  # tx30 y 3928099 -> 3883643
$grp.GroupItems.Item(29).Top = 221.2128448486328  # tx31 y 2876088 -> 2809403
$grp.GroupItems.Item(30).Top = 136.638671875  # tx32 y 1824225 -> 1735311
$grp.GroupItems.Item(32).Top = 394.855224609375  # pl34 y 5036889 -> 5014661
$grp.GroupItems.Item(33).Top = 310.27520751953125  # pl35 y 3984952 -> 3940495
$grp.GroupItems.Item(34).Top = 225.69520568847656  # pl36 y 2933014 -> 2866329
$grp.GroupItems.Item(35).Top = 141.1152801513672  # pl37 y 1881077 -> 1792164
$grp.GroupItems.Item(51).TextFrame.TextRange.Text = "(N=143)"  # tx53 text '(N=146)' -> '(N=143)'
